# Insert two new rows of "Ciruela" price data (Black Amber, O'Higgins, 2022-01-25)
# above the existing row 79, pushing the old rows 79-164 down to 81-166.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 79 (existing rows 79:164 shift down to 81:166)
$ws.Rows("79:80").Insert()

# --- New row 79 ---
$ws.Range("A79").Value = 8
$ws.Range("B79").Value = "Terminal La Palmera de La Serena"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44586
$ws.Range("E79").Value = 4
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100103
$ws.Range("H79").Value = "Frutos de hueso (carozo)"
$ws.Range("I79").Value = 100103002
$ws.Range("J79").Value = "Ciruela"
$ws.Range("K79").Value = "Black Amber"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 24
$ws.Range("N79").Value = 260000
$ws.Range("O79").Value = 270000
$ws.Range("P79").Value = 265000
$ws.Range("Q79").Value = "`$/bins (450 kilos)"
$ws.Range("R79").Value = "Región de O'Higgins"
$ws.Range("S79").Value = 589
$ws.Range("T79").Value = 450

# --- New row 80 ---
$ws.Range("A80").Value = 8
$ws.Range("B80").Value = "Terminal La Palmera de La Serena"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44586
$ws.Range("E80").Value = 4
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100103
$ws.Range("H80").Value = "Frutos de hueso (carozo)"
$ws.Range("I80").Value = 100103002
$ws.Range("J80").Value = "Ciruela"
$ws.Range("K80").Value = "Black Amber"
$ws.Range("L80").Value = "Segunda"
$ws.Range("M80").Value = 20
$ws.Range("N80").Value = 210000
$ws.Range("O80").Value = 220000
$ws.Range("P80").Value = 215000
$ws.Range("Q80").Value = "`$/bins (450 kilos)"
$ws.Range("R80").Value = "Región de O'Higgins"
$ws.Range("S80").Value = 478
$ws.Range("T80").Value = 450
